$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.257.18'
$ws.Range("E2").Value = '  -1.13%  '
$ws.Range("D3").Value = '2.234.78'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'229.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("E6").Value = '  +1.54%  '
$ws.Range("D7").Value = "'63.60"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.59%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = "'0.436"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("E10").Value = '  -5.80%  '
$ws.Range("D11").Value = "'56.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("D12").Value = "'26.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.76%  '
$ws.Range("E13").Value = '  -1.71%  '
$ws.Range("D14").Value = '2.570.32'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").Value = "'15.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.13%  '
$ws.Range("D16").Value = "'5.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = "'0.818"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.26%  '
$ws.Range("D18").Value = '2.235.11'
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").Value = '43.178.74'
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("D20").Value = '0.0₃0956'
$ws.Range("E20").Value = '  -5.66%  '
$ws.Range("D21").Value = "'72.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").Value = "'6.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = "'245.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.83%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = "'3.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +29.95%  '
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("E27").Value = '  -1.78%  '
$ws.Range("D28").Value = "'173.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.58%  '
$ws.Range("D29").Value = "'9.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.37%  '
$ws.Range("D30").Value = "'21.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.29%  '
$ws.Range("D31").Value = "'0.128"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.31%  '
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("D34").Value = "'4.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.11%  '
$ws.Range("D35").Value = "'0.0671"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.34%  '
$ws.Range("D36").Value = "'4.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.99%  '
$ws.Range("E37").Value = '  -7.33%  '
$ws.Range("D38").Value = "'6.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.57%  '
$ws.Range("D39").Value = "'2.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.92%  '
$ws.Range("D40").Value = "'0.0247"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = "'8.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.73%  '
$ws.Range("D43").Value = "'4.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.74%  '
$ws.Range("D44").Value = "'16.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.69%  '
$ws.Range("D45").Value = "'95.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.84%  '
$ws.Range("D46").Value = "'0.0935"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.37%  '
$ws.Range("D47").Value = "'1.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("E48").Value = '  -2.39%  '
$ws.Range("D49").Value = '1.420.05'
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("D50").Value = "'9.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("E51").Value = '  +0.41%  '
